$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the mismatched F41 value (event "An Absence of Mind 2") ---
$ws.Range("F41").Value = -2

# --- Append the new random-event rows (46-55) with their stat deltas ---
$newEvents = @(
    @{ Row = 46; Name = "A Scurrying Racket 1"; F = 0;  G = 0;  H = -1 },
    @{ Row = 47; Name = "A Scurrying Racket 2"; F = -2; G = 0;  H = 2 },
    @{ Row = 48; Name = "Deepening Chill 1";    F = 1;  G = 0;  H = 1 },
    @{ Row = 49; Name = "Deepening Chill 2";    F = 2;  G = -1; H = 1 },
    @{ Row = 50; Name = "Unknown Chuckle 1";    F = 2;  G = -2; H = 0 },
    @{ Row = 51; Name = "Unknown Chuckle 2";    F = 0;  G = 2;  H = -2 },
    @{ Row = 52; Name = "Withered Petals 1";    F = 1;  G = -1; H = 0 },
    @{ Row = 53; Name = "Withered Petals 2";    F = 0;  G = 1;  H = -1 },
    @{ Row = 54; Name = "Liquid Yearning 1";    F = -3; G = 2;  H = 0 },
    @{ Row = 55; Name = "Liquid Yearning 2";    F = 0;  G = -2; H = 1 }
)

foreach ($evt in $newEvents) {
    $r = $evt.Row
    $ws.Range("D$r").Value = $evt.Name
    $ws.Range("F$r").Value = $evt.F
    $ws.Range("G$r").Value = $evt.G
    $ws.Range("H$r").Value = $evt.H
}

# --- Extend the averaging formulas to cover the newly added rows ---
$ws.Range("J26").Formula = "=AVERAGE(F26:F55)"
$ws.Range("K26").Formula = "=AVERAGE(G26:G55)"
$ws.Range("L26").Formula = "=AVERAGE(H26:H55)"

# --- Update the view: selected cell and top-left visible cell ---
$ws.Range("M30").Select()
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 4

$excel.Calculate()
